$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 3")

# Row 10 additions
$ws.Range("D10").Value = 0.8979166666666667
$ws.Range("E10").Value = 90
$ws.Range("F10").Value = 408
$ws.Range("H10").Value = "p. 3 - 7"
$ws.Range("J10").Value = "x"

# Row 11 additions
$ws.Range("B11").Value = 43877
$ws.Range("C11").Value = 0.98263888888888884

# Update selection to match the diff
$ws.Range("B11").Select()
